$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 545, shifting existing rows 545:587 down to 546:588.
$ws.Rows(545).Insert()

# Populate the new row 545 with the new weekly price record. The fields that
# are identical across the whole dataset (A,B,C,E,F,G,H,I,N,Q,R) are copied
# from the neighbouring row; only D,J,K,L,M,O,P carry new data per the diff.
$ws.Range("A545").Value = 4
$ws.Range("B545").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C545").Value = "Los Lagos"
$ws.Range("D545").Value = 45106
$ws.Range("E545").Value = 10
$ws.Range("F545").Value = 100114013
$ws.Range("G545").Value = "Zanahoria"
$ws.Range("H545").Value = "Sin especificar"
$ws.Range("I545").Value = "Primera"
$ws.Range("J545").Value = 150
$ws.Range("K545").Value = 9000
$ws.Range("L545").Value = 9000
$ws.Range("M545").Value = 9000
$ws.Range("N545").Value = "$/saco 20 kilos"
$ws.Range("O545").Value = "Provincia de Llanquihue"
$ws.Range("P545").Value = 450
$ws.Range("Q545").Value = 20
$ws.Range("R545").Value = "Hortaliza"
